# Edit script: update the cached "datetimeFigureOut" date field text in the
# slide master and all slide layouts (06.08.2025 -> 23.09.2025), and
# reposition/resize four shapes on slide 2 (the WBZ451H modifications slide)
# to make room for / reflect the newly-included Adafruit protomatter library
# annotation.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached date text on the slide master's Date Placeholder.
# ---------------------------------------------------------------------
$m = $p.SlideMaster

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

Update-DatePlaceholder $m.Shapes "23.09.2025"

# ---------------------------------------------------------------------
# 2) Update the cached date text on every slide layout's Date Placeholder.
# ---------------------------------------------------------------------
for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $cl = $m.CustomLayouts.Item($j)
    Update-DatePlaceholder $cl.Shapes "23.09.2025"
}

# ---------------------------------------------------------------------
# 3) Reposition / resize shapes on slide 2.
#
# Shape Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU)
# and are stored internally as single-precision floats that get
# truncated (not rounded) back to EMU, so a plain "emu / 12700" can
# land 1 EMU short. Adding half an EMU's worth of points before
# dividing compensates for that truncation and reproduces the exact
# EMU values from the target OOXML.
# ---------------------------------------------------------------------
$EMU_PER_POINT = 12700

function EmuToPt($emu) {
    return ($emu + 0.5) / $EMU_PER_POINT
}

$s2 = $p.Slides.Item(2)

# "Rectangle 3" (green highlight box)
$rect3 = $s2.Shapes.Item(2)
$rect3.Left = EmuToPt 6585286
$rect3.Top = EmuToPt 3360821
$rect3.Width = EmuToPt 751686
$rect3.Height = EmuToPt 705853

# "Straight Arrow Connector 7"
$conn7 = $s2.Shapes.Item(4)
$conn7.Left = EmuToPt 2845183
$conn7.Top = EmuToPt 4143680
$conn7.Width = EmuToPt 715887
$conn7.Height = EmuToPt 348109

# "Rectangle 19" (green highlight box)
$rect19 = $s2.Shapes.Item(11)
$rect19.Left = EmuToPt 3305534
$rect19.Top = EmuToPt 4950259
$rect19.Width = EmuToPt 320841
$rect19.Height = EmuToPt 260684

# "Straight Arrow Connector 20" - also gains a horizontal flip
$conn20 = $s2.Shapes.Item(12)
$conn20.Left = EmuToPt 3626375
$conn20.Top = EmuToPt 5210943
$conn20.Width = EmuToPt 221439
$conn20.Height = EmuToPt 911626
$conn20.HorizontalFlip = $true
